$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "level": clear old level-layout data (columns A-D, rows 1-23) and
# write the new level layout ("levels up to 50"). Column E keeps its
# existing formula (=meta!$B$1*ROW()) and is left untouched.
# ---------------------------------------------------------------------------
$wsLevel = $wb.Worksheets.Item("level")
$wsLevel.Range("A1:D23").ClearContents()

$wsLevel.Range("A1").Value = 7
$wsLevel.Range("D1").Value = 7

$wsLevel.Range("A3").Value = 7
$wsLevel.Range("D3").Value = 8

$wsLevel.Range("B4").Value = 8
$wsLevel.Range("C4").Value = 7

$wsLevel.Range("B6").Value = 10

$wsLevel.Range("A7").Value = 8
$wsLevel.Range("B7").Value = 7
$wsLevel.Range("C7").Value = 7
$wsLevel.Range("D7").Value = 8

$wsLevel.Range("C9").Value = 11

$wsLevel.Range("A10").Value = 7
$wsLevel.Range("B10").Value = 8
$wsLevel.Range("C10").Value = 7
$wsLevel.Range("D10").Value = 7

# ---------------------------------------------------------------------------
# Sheet "enemies": rework the enemy stat table (rows 2-8) and append new
# enemy entries (rows 9-12) including the new "Butterfly"-type enemies and
# the new G-column resource paths for rows 6-7.
# ---------------------------------------------------------------------------
$wsEnemies = $wb.Worksheets.Item("enemies")

$wsEnemies.Range("B2").Value = 0
$wsEnemies.Range("C2").Value = 0
$wsEnemies.Range("D2").Value = 2
$wsEnemies.Range("E2").Value = 1
$wsEnemies.Range("F2").Value = "EnemyPrefabs/Special Enemies/Halloween Bee/Halloween Bee"

$wsEnemies.Range("B3").Value = 0
$wsEnemies.Range("C3").Value = 0
$wsEnemies.Range("D3").Value = 2
$wsEnemies.Range("E3").Value = 1
$wsEnemies.Range("F3").Value = "EnemyPrefabs/Special Enemies/Steampunk Fly/Steampunk Fly"

$wsEnemies.Range("B4").Value = 0
$wsEnemies.Range("C4").Value = 0
$wsEnemies.Range("D4").Value = 2
$wsEnemies.Range("E4").Value = 1
$wsEnemies.Range("F4").Value = "EnemyPrefabs/Special Enemies/Bionic Lady Bird/Bionic Lady Bird"

$wsEnemies.Range("B5").Value = 0
$wsEnemies.Range("C5").Value = 0
$wsEnemies.Range("D5").Value = 3
$wsEnemies.Range("E5").Value = 1.5
$wsEnemies.Range("F5").Value = "EnemyPrefabs/Special Enemies/Halloween Bee/Halloween Bee"

$wsEnemies.Range("B6").Value = 0
$wsEnemies.Range("C6").Value = 0
$wsEnemies.Range("D6").Value = 3
$wsEnemies.Range("E6").Value = 1.5
$wsEnemies.Range("F6").Value = "EnemyPrefabs/Special Enemies/Steampunk Fly/Steampunk Fly"
$wsEnemies.Range("G6").Value = "EnemyPrefabs/Bullet Enemies//"

$wsEnemies.Range("B7").Value = 0
$wsEnemies.Range("C7").Value = 0
$wsEnemies.Range("D7").Value = 3
$wsEnemies.Range("E7").Value = 1.5
$wsEnemies.Range("F7").Value = "EnemyPrefabs/Special Enemies/Bionic Lady Bird/Bionic Lady Bird"
$wsEnemies.Range("G7").Value = "EnemyPrefabs/Special Enemies//"

$wsEnemies.Range("B8").Value = 4
$wsEnemies.Range("C8").Value = 0
$wsEnemies.Range("D8").Value = 0
$wsEnemies.Range("E8").Value = 1
$wsEnemies.Range("F8").Value = "EnemyPrefabs/Arrow Enemies/Bee/Bee Arrow"

$wsEnemies.Range("A9").Value = 8
$wsEnemies.Range("B9").Value = 3
$wsEnemies.Range("C9").Value = 0
$wsEnemies.Range("D9").Value = 0
$wsEnemies.Range("E9").Value = 1.2
$wsEnemies.Range("F9").Value = "EnemyPrefabs/Arrow Enemies/Butterfly/Butterfly Arrow"

$wsEnemies.Range("A10").Value = 9
$wsEnemies.Range("B10").Value = 3
$wsEnemies.Range("C10").Value = 0
$wsEnemies.Range("D10").Value = 0
$wsEnemies.Range("E10").Value = 1.2
$wsEnemies.Range("F10").Value = "EnemyPrefabs/Bullet Enemies/Purple Butterfly/Butterfly Arrow"

$wsEnemies.Range("A11").Value = 10
$wsEnemies.Range("B11").Value = 3
$wsEnemies.Range("C11").Value = 0
$wsEnemies.Range("D11").Value = 0
$wsEnemies.Range("E11").Value = 1.2
$wsEnemies.Range("F11").Value = "EnemyPrefabs/Bullet Enemies/Purple Butterfly/Purple Butterfly"

$wsEnemies.Range("A12").Value = 11
$wsEnemies.Range("B12").Value = 3
$wsEnemies.Range("C12").Value = 0
$wsEnemies.Range("D12").Value = 0
$wsEnemies.Range("E12").Value = 1.2
$wsEnemies.Range("F12").Value = "EnemyPrefabs/Special Enemies/Red Butterfly/Red Butterfly"

# ---------------------------------------------------------------------------
# Sheet "misc": fill in the previously empty row 5 with the new Butterfly
# Arrow resource path (list of "all enemies").
# ---------------------------------------------------------------------------
$wsMisc = $wb.Worksheets.Item("misc")
$wsMisc.Range("A5").Value = "EnemyPrefabs/Arrow Enemies/Butterfly/Butterfly Arrow"
